$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 105, shifting rows 105:204 down to 106:205
$ws.Rows("105:105").Insert()

# Populate new row 105 with data (copy of the static columns from row 106, new D/J values)
$ws.Range("A105").Value = 4
$ws.Range("B105").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C105").Value = "Los Lagos"
$ws.Range("D105").Value = 44587
$ws.Range("E105").Value = 10
$ws.Range("F105").Value = 100112044
$ws.Range("G105").Value = "Perejil"
$ws.Range("H105").Value = "Sin especificar"
$ws.Range("I105").Value = "Primera"
$ws.Range("J105").Value = 20
$ws.Range("K105").Value = 6000
$ws.Range("L105").Value = 6000
$ws.Range("M105").Value = 6000
$ws.Range("N105").Value = "$/docena de atados (2 kilos)"
$ws.Range("O105").Value = "Región de La Araucanía"
$ws.Range("P105").Value = 3000
$ws.Range("Q105").Value = 2
$ws.Range("R105").Value = "Hortaliza"

# Ensure D105 uses the same date style as other D cells (style copied from insert, but set explicitly just in case)
$ws.Range("D105").NumberFormat = $ws.Range("D106").NumberFormat
